$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: Intel(R) Dual Band Wireless-AC 7260 - 17.15.0.5 -- Critical Minutes updated
$ws.Range("C5").Value = 267

# Row 6: now holds the "7265 - 19.51.42.2" adapter with updated counts
$ws.Range("A6").Value = "Intel(R) Dual Band Wireless-AC 7265 - 19.51.42.2"
$ws.Range("B6").Value = 36
$ws.Range("C6").Value = 2558

# Row 7: now holds the "8265 - 20.70.32.1" adapter with updated counts
$ws.Range("A7").Value = "Intel(R) Dual Band Wireless-AC 8265 - 20.70.32.1"
$ws.Range("B7").Value = 8
$ws.Range("C7").Value = 550

# Row 8: Totals updated to reflect new sums
$ws.Range("B8").Value = 53
$ws.Range("C8").Value = 6437
